$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.775.29"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").Value = "'2.577.11"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -1.78%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = "'584.38"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").Value = "'168.89"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +1.50%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("E8").Value = '  -0.80%  '

$ws.Range("D9").Value = "'2.576.57"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -1.85%  '

$ws.Range("E10").Value = '  +0.37%  '

$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").Value = "'5.17"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("D14").Value = "'26.86"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").Value = "'3.048.49"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.96%  '

$ws.Range("E16").Value = '  -1.32%  '

$ws.Range("D17").Value = "'66.544.93"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").Value = "'2.584.26"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -1.56%  '

$ws.Range("D19").Value = "'11.44"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -5.69%  '

$ws.Range("D20").Value = "'7.76"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -3.35%  '

$ws.Range("D21").Value = "'351.28"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.28%  '

$ws.Range("D22").Value = "'4.24"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").Value = "'4.62"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").Value = "'1.92"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("D26").Value = "'69.41"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("D27").Value = "'9.93"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -8.80%  '

$ws.Range("E28").Value = '  -1.88%  '

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("D30").Value = "'0.0₃0994"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -0.94%  '

$ws.Range("D31").Value = "'532.61"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -2.78%  '

$ws.Range("D32").Value = "'8.26"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +4.46%  '

$ws.Range("E33").Value = '  -1.56%  '

$ws.Range("E34").Value = '  -3.10%  '

$ws.Range("E35").Value = '  -2.74%  '

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("D37").Value = "'1.47"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -2.01%  '

$ws.Range("D38").Value = "'156.79"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("E39").Value = '  -1.15%  '

$ws.Range("E40").Value = '  -1.85%  '

$ws.Range("D41").Value = "'18.34"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +2.23%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.78"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = "'5.14"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -0.28%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = "'2.44"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.69%  '

$ws.Range("D46").Value = "'0.0₆0287"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -2.92%  '

$ws.Range("D47").Value = "'149.54"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -1.11%  '

$ws.Range("D48").Value = "'0.568"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("D49").Value = "'3.73"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("E50").Value = '  +1.07%  '

$ws.Range("E51").Value = '  -0.87%  '
